$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append bibliography rows 25-35 (new shared strings 23-33) ---

$ws.Range("A25").Value2 = "Rehm, H. L. (2017). Evolving health care through personal genomics. Nature Reviews Genetics."
$ws.Range("A25").Font.Name = "Arial"
$ws.Range("A25").Font.Size = 8
$ws.Range("A25").Font.Color = 2236962
$chars_25_1 = $ws.Range("A25").Characters(69, 23)
$chars_25_1.Font.Italic = $true
$chars_25_1.Font.Size = 8
$chars_25_1.Font.Name = "Arial"
$chars_25_1.Font.Color = 2236962
$chars_25_2 = $ws.Range("A25").Characters(92, 1)
$chars_25_2.Font.Size = 8
$chars_25_2.Font.Name = "Arial"
$chars_25_2.Font.Color = 2236962

$ws.Range("A26").Value2 = "Angrist, M. (2016). Personal genomics: Where are we now?. Applied & translational genomics, 8, 1."
$ws.Range("A26").Font.Name = "Arial"
$ws.Range("A26").Font.Size = 8
$ws.Range("A26").Font.Color = 2236962
$chars_26_1 = $ws.Range("A26").Characters(59, 32)
$chars_26_1.Font.Italic = $true
$chars_26_1.Font.Size = 8
$chars_26_1.Font.Name = "Arial"
$chars_26_1.Font.Color = 2236962
$chars_26_2 = $ws.Range("A26").Characters(91, 2)
$chars_26_2.Font.Size = 8
$chars_26_2.Font.Name = "Arial"
$chars_26_2.Font.Color = 2236962
$chars_26_3 = $ws.Range("A26").Characters(93, 1)
$chars_26_3.Font.Italic = $true
$chars_26_3.Font.Size = 8
$chars_26_3.Font.Name = "Arial"
$chars_26_3.Font.Color = 2236962
$chars_26_4 = $ws.Range("A26").Characters(94, 4)
$chars_26_4.Font.Size = 8
$chars_26_4.Font.Name = "Arial"
$chars_26_4.Font.Color = 2236962

$ws.Range("A27").Value2 = "Escalona, M., Rocha, S., & Posada, D. (2016). A comparison of tools for the simulation of genomic next-generation sequencing data. Nature Reviews Genetics, 17(8), 459-469."
$ws.Range("A27").Font.Name = "Arial"
$ws.Range("A27").Font.Size = 8
$ws.Range("A27").Font.Color = 2236962
$chars_27_1 = $ws.Range("A27").Characters(132, 23)
$chars_27_1.Font.Italic = $true
$chars_27_1.Font.Size = 8
$chars_27_1.Font.Name = "Arial"
$chars_27_1.Font.Color = 2236962
$chars_27_2 = $ws.Range("A27").Characters(155, 2)
$chars_27_2.Font.Size = 8
$chars_27_2.Font.Name = "Arial"
$chars_27_2.Font.Color = 2236962
$chars_27_3 = $ws.Range("A27").Characters(157, 2)
$chars_27_3.Font.Italic = $true
$chars_27_3.Font.Size = 8
$chars_27_3.Font.Name = "Arial"
$chars_27_3.Font.Color = 2236962
$chars_27_4 = $ws.Range("A27").Characters(159, 13)
$chars_27_4.Font.Size = 8
$chars_27_4.Font.Name = "Arial"
$chars_27_4.Font.Color = 2236962

$ws.Range("A28").Value2 = "Schirmer, M., D’Amore, R., Ijaz, U. Z., Hall, N., & Quince, C. (2016). Illumina error profiles: resolving fine-scale variation in metagenomic sequencing data. BMC bioinformatics, 17(1), 125."
$ws.Range("A28").Font.Name = "Arial"
$ws.Range("A28").Font.Size = 8
$ws.Range("A28").Font.Color = 2236962
$chars_28_1 = $ws.Range("A28").Characters(160, 18)
$chars_28_1.Font.Italic = $true
$chars_28_1.Font.Size = 8
$chars_28_1.Font.Name = "Arial"
$chars_28_1.Font.Color = 2236962
$chars_28_2 = $ws.Range("A28").Characters(178, 2)
$chars_28_2.Font.Size = 8
$chars_28_2.Font.Name = "Arial"
$chars_28_2.Font.Color = 2236962
$chars_28_3 = $ws.Range("A28").Characters(180, 2)
$chars_28_3.Font.Italic = $true
$chars_28_3.Font.Size = 8
$chars_28_3.Font.Name = "Arial"
$chars_28_3.Font.Color = 2236962
$chars_28_4 = $ws.Range("A28").Characters(182, 9)
$chars_28_4.Font.Size = 8
$chars_28_4.Font.Name = "Arial"
$chars_28_4.Font.Color = 2236962

$ws.Range("A29").Value2 = "LeCun, Y., Bengio, Y., & Hinton, G. (2015). Deep learning. Nature, 521(7553), 436-444."
$ws.Range("A29").Font.Name = "Arial"
$ws.Range("A29").Font.Size = 8
$ws.Range("A29").Font.Color = 2236962
$chars_29_1 = $ws.Range("A29").Characters(60, 6)
$chars_29_1.Font.Italic = $true
$chars_29_1.Font.Size = 8
$chars_29_1.Font.Name = "Arial"
$chars_29_1.Font.Color = 2236962
$chars_29_2 = $ws.Range("A29").Characters(66, 2)
$chars_29_2.Font.Size = 8
$chars_29_2.Font.Name = "Arial"
$chars_29_2.Font.Color = 2236962
$chars_29_3 = $ws.Range("A29").Characters(68, 3)
$chars_29_3.Font.Italic = $true
$chars_29_3.Font.Size = 8
$chars_29_3.Font.Name = "Arial"
$chars_29_3.Font.Color = 2236962
$chars_29_4 = $ws.Range("A29").Characters(71, 16)
$chars_29_4.Font.Size = 8
$chars_29_4.Font.Name = "Arial"
$chars_29_4.Font.Color = 2236962

$ws.Range("A30").Value2 = "Kingma, D., & Ba, J. (2014). Adam: A method for stochastic optimization. arXiv preprint arXiv:1412.6980."
$ws.Range("A30").Font.Name = "Arial"
$ws.Range("A30").Font.Size = 8
$ws.Range("A30").Font.Color = 2236962
$chars_30_1 = $ws.Range("A30").Characters(74, 30)
$chars_30_1.Font.Italic = $true
$chars_30_1.Font.Size = 8
$chars_30_1.Font.Name = "Arial"
$chars_30_1.Font.Color = 2236962
$chars_30_2 = $ws.Range("A30").Characters(104, 1)
$chars_30_2.Font.Size = 8
$chars_30_2.Font.Name = "Arial"
$chars_30_2.Font.Color = 2236962

$ws.Range("A31").Value2 = "Sutskever, I., Martens, J., Dahl, G. E., & Hinton, G. E. (2013). On the importance of initialization and momentum in deep learning. ICML (3), 28, 1139-1147."
$ws.Range("A31").Font.Name = "Arial"
$ws.Range("A31").Font.Size = 8
$ws.Range("A31").Font.Color = 2236962

$ws.Range("A32").Value2 = "Maas, A. L., Hannun, A. Y., & Ng, A. Y. (2013, June). Rectifier nonlinearities improve neural network acoustic models. In Proc. ICML (Vol. 30, No. 1)."
$ws.Range("A32").Font.Name = "Arial"
$ws.Range("A32").Font.Size = 8
$ws.Range("A32").Font.Color = 2236962
$chars_32_1 = $ws.Range("A32").Characters(123, 10)
$chars_32_1.Font.Italic = $true
$chars_32_1.Font.Size = 8
$chars_32_1.Font.Name = "Arial"
$chars_32_1.Font.Color = 2236962
$chars_32_2 = $ws.Range("A32").Characters(133, 18)
$chars_32_2.Font.Size = 8
$chars_32_2.Font.Name = "Arial"
$chars_32_2.Font.Color = 2236962

$ws.Range("A33").Value2 = "Srivastava, N., Hinton, G. E., Krizhevsky, A., Sutskever, I., & Salakhutdinov, R. (2014). Dropout: a simple way to prevent neural networks from overfitting. Journal of Machine Learning Research, 15(1), 1929-1958."
$ws.Range("A33").Font.Name = "Arial"
$ws.Range("A33").Font.Size = 8
$ws.Range("A33").Font.Color = 2236962

$ws.Range("A34").Value2 = "Ruder, S. (2016). An overview of gradient descent optimization algorithms. arXiv preprint arXiv:1609.04747."
$ws.Range("A34").Font.Name = "Arial"
$ws.Range("A34").Font.Size = 8
$ws.Range("A34").Font.Color = 2236962
$chars_34_1 = $ws.Range("A34").Characters(76, 31)
$chars_34_1.Font.Italic = $true
$chars_34_1.Font.Size = 8
$chars_34_1.Font.Name = "Arial"
$chars_34_1.Font.Color = 2236962
$chars_34_2 = $ws.Range("A34").Characters(107, 1)
$chars_34_2.Font.Size = 8
$chars_34_2.Font.Name = "Arial"
$chars_34_2.Font.Color = 2236962

$ws.Range("A35").Value2 = "Tieleman, T. and Hinton, G. Lecture 6.5 - RMSProp, COURSERA: Neural Networks for Machine Learning.`nTechnical report, 2012"

# --- Update view: selection moves to K28 ---
$ws.Range("K28").Select()
